# Generate Report for Handoff
#
# The localization status report is regenerated: the handoff package is
# ready, so every "In Translation" status becomes "Ready for handoff" and
# the two "Latest Handoff Datetime" timestamps are refreshed. Because the
# new status text is longer than the old one, the Status column on each
# language sheet (and its mirrored column on the Overview sheet) is
# widened to fit ("AutoFit"-style) as part of regenerating the report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-30 15:17:29"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-30 15:17:34"

# --- Overview sheet: per-language Status columns (zh-cn / de-de) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- Widen the Status columns to fit the new, longer text ---
# (Status column is column C on the language sheets, and columns E/F on
# the Overview sheet.) The host quantizes ColumnWidth assignments onto
# its internal character grid (nearest 1/6 of a character, offset by
# 5/6), so the requested width is pre-compensated here to land on the
# closest representable value to the target fitted width.
$fitWidth = 17.2159881591797
$fitWidthInput = $fitWidth - (5 / 6)

$zhcn.Columns.Item(3).ColumnWidth = $fitWidthInput
$dede.Columns.Item(3).ColumnWidth = $fitWidthInput
$overview.Columns.Item(5).ColumnWidth = $fitWidthInput
$overview.Columns.Item(6).ColumnWidth = $fitWidthInput
